$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "crafting section" parent references (rows 31-33) ---
$ws.Range("P31").Value = "The Return of the King's Crown"
$ws.Range("P32").Value = "The Return of the King's Crown"
$ws.Range("P33").Value = "Dead Animals"

# --- Add new quest rows (52-66) for the new skill quest items ---

# Row 52
$ws.Range("A52").Value = "Gobbies and Kings"
$ws.Range("B52").Value = "Helpless Goblin"
$ws.Range("D52").Value = 1000
$ws.Range("H52").Value = "Kings Book of Hope"

# Row 53
$ws.Range("A53").Value = "Move with the wind"
$ws.Range("B53").Value = "Helpless Goblin"
$ws.Range("D53").Value = 1000
$ws.Range("H53").Value = "Feathers of the sky"
$ws.Range("K53").Value = 100
$ws.Range("L53").Value = 5
$ws.Range("P53").Value = "Goblins Lust for Gold"

# Row 54
$ws.Range("A54").Value = "Story of the wind"
$ws.Range("B54").Value = "Helpless Goblin"
$ws.Range("C54").Value = "Feathers of the sky"
$ws.Range("D54").Value = 1500
$ws.Range("H54").Value = "Book of speed"
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 10
$ws.Range("P54").Value = "Goblins Lust for Gold"
$ws.Range("R54").Value = "Surface"
$ws.Range("S54").Value = 1

# Row 55
$ws.Range("A55").Value = "Sky is where the dragons live"
$ws.Range("B55").Value = "Helpless Goblin"
$ws.Range("C55").Value = "Book of speed"
$ws.Range("D55").Value = 3000
$ws.Range("H55").Value = "Sky Dragons Blood Vial"
$ws.Range("K55").Value = 100000
$ws.Range("L55").Value = 50
$ws.Range("P55").Value = "Story of the wind"

# Row 56
$ws.Range("A56").Value = "The Wizards Enchantment"
$ws.Range("B56").Value = "Helpless Goblin"
$ws.Range("C56").Value = "Sky Dragons Blood Vial"
$ws.Range("D56").Value = 6000
$ws.Range("H56").Value = "The Wizards Enchantment"
$ws.Range("K56").Value = 1000000
$ws.Range("L56").Value = 100
$ws.Range("P56").Value = "Sky is where the dragons live"

# Row 57
$ws.Range("A57").Value = "Lost Arrow"
$ws.Range("B57").Value = "Helpless Goblin"
$ws.Range("C57").Value = "Heart of the Eye"
$ws.Range("D57").Value = 1000
$ws.Range("H57").Value = "Arrow of Truth"
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 10
$ws.Range("P57").Value = "Goblins Lust for Gold"
$ws.Range("R57").Value = "Surface"
$ws.Range("S57").Value = 3

# Row 58
$ws.Range("A58").Value = "Goblins Goddess"
$ws.Range("B58").Value = "Helpless Goblin"
$ws.Range("C58").Value = "Arrow of Truth"
$ws.Range("D58").Value = 1500
$ws.Range("H58").Value = "Goddess Bow of Light"
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 50
$ws.Range("P58").Value = "Lost Arrow"

# Row 59
$ws.Range("A59").Value = "Goblins Accuracy"
$ws.Range("B59").Value = "Helpless Goblin"
$ws.Range("C59").Value = "Goddess Bow of Light"
$ws.Range("D59").Value = 3000
$ws.Range("H59").Value = "The Book of Guidance"
$ws.Range("K59").Value = 1000000
$ws.Range("L59").Value = 100
$ws.Range("P59").Value = "Goblins Goddess"

# Row 60
$ws.Range("A60").Value = "Kings and Queens"
$ws.Range("B60").Value = "Helpless Goblin"
$ws.Range("C60").Value = "Kings Book of Hope"
$ws.Range("D60").Value = 1500
$ws.Range("H60").Value = "Settlers Walking Stick"
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 10
$ws.Range("P60").Value = "Goblins Lust for Gold"
$ws.Range("R60").Value = "Surface"
$ws.Range("S60").Value = 5

# Row 61
$ws.Range("A61").Value = "Take a walk"
$ws.Range("B61").Value = "Helpless Goblin"
$ws.Range("C61").Value = "Settlers Walking Stick"
$ws.Range("D61").Value = 3000
$ws.Range("H61").Value = "Kings Scepter"
$ws.Range("K61").Value = 100000
$ws.Range("L61").Value = 50
$ws.Range("P61").Value = "Kings and Queens"

# Row 62
$ws.Range("A62").Value = "Leading Armies"
$ws.Range("B62").Value = "Helpless Goblin"
$ws.Range("C62").Value = "Kings Scepter"
$ws.Range("D62").Value = 5000
$ws.Range("H62").Value = "Kings Ring"
$ws.Range("K62").Value = 1000000
$ws.Range("L62").Value = 100
$ws.Range("P62").Value = "Take a walk"

# Row 63
$ws.Range("A63").Value = "Alchemist Scrolls"
$ws.Range("B63").Value = "The Witch"
$ws.Range("C63").Value = "Alchemist Book"
$ws.Range("D63").Value = 5000
$ws.Range("E63").Value = 5
$ws.Range("H63").Value = "Alchemist Scrolls"
$ws.Range("I63").Value = 10000
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 50
$ws.Range("O63").Value = 1
$ws.Range("R63").Value = "Labyrinth"
$ws.Range("S63").Value = 2

# Row 64
$ws.Range("A64").Value = "Run forest, run!"
$ws.Range("B64").Value = "Helpless Goblin"
$ws.Range("C64").Value = "Feathers of the sky"
$ws.Range("D64").Value = 1000
$ws.Range("H64").Value = "Mages Teleport Scroll"
$ws.Range("K64").Value = 1000
$ws.Range("L64").Value = 10
$ws.Range("P64").Value = "Move with the wind"
$ws.Range("R64").Value = "Labyrinth"
$ws.Range("S64").Value = 3

# Row 65
$ws.Range("A65").Value = "Blink of an eye"
$ws.Range("B65").Value = "Helpless Goblin"
$ws.Range("C65").Value = "Mages Teleport Scroll"
$ws.Range("D65").Value = 2000
$ws.Range("E65").Value = 10
$ws.Range("H65").Value = "Goblins Quickening Rune"
$ws.Range("J65").Value = 20
$ws.Range("K65").Value = 100000
$ws.Range("L65").Value = 50
$ws.Range("P65").Value = "Run forest, run!"

# Row 66
$ws.Range("A66").Value = "Goblins Gift"
$ws.Range("B66").Value = "Helpless Goblin"
$ws.Range("C66").Value = "Goblins Quickening Rune"
$ws.Range("D66").Value = 5000
$ws.Range("E66").Value = 10
$ws.Range("H66").Value = "Golden Ring of Blur"
$ws.Range("I66").Value = 10000
$ws.Range("J66").Value = 100
$ws.Range("K66").Value = 1000000
$ws.Range("L66").Value = 100
$ws.Range("P66").Value = "Blink of an eye"
$ws.Range("R66").Value = "Dungeons"
$ws.Range("S66").Value = 1
